# Simulated Wild Card round and logged it
# Update the Target Depth Data for the Cowboys: OFF ("H" row) and DEF ("H" row)
# reflect the added Wild Card game's passing-depth stats.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 486   # Short Att
$wsOff.Range("C2").Value = 330   # Short Comp
$wsOff.Range("D2").Value = 118   # Deep Att
$wsOff.Range("E2").Value = 70    # Deep Comp
$wsOff.Range("F2").Value = 8     # Short Int

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 323   # Short Att
$wsDef.Range("C2").Value = 215   # Short Comp
$wsDef.Range("D2").Value = 80    # Deep Att
$wsDef.Range("E2").Value = 44    # Deep Comp
$wsDef.Range("F2").Value = 8     # Short Int

$wb.Save()
